# Update election result values for BRAGANÇA / BRAGANÇA row (row 2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "H2"  = 398
    "I2"  = 1052
    "J2"  = 4396
    "K2"  = 24
    "L2"  = 1194
    "M2"  = 72
    "N2"  = 802
    "O2"  = 3
    "P2"  = 17
    "Q2"  = 8
    "R2"  = 50
    "S2"  = 505
    "T2"  = 759
    "U2"  = 54
    "V2"  = 6718
    "W2"  = 5
    "X2"  = 6716
    "Y2"  = 12
    "Z2"  = 117
    "AA2" = 51
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
